$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows (participants 6 and 7) appended to the response table ---

# Row 7 (participant ID 6)
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 44573.701666666697
$ws.Cells.Item(7, 3).Value = 44573.706921296303
$ws.Cells.Item(7, 4).Value = "'6"
$ws.Cells.Item(7, 5).Value = "Pinch Anywhere;Dwell;Pinch on Circle;Touch In The Air;"
$ws.Cells.Item(7, 6).Value = "This interaction was hard to use and the least natural and intuitive method. Its awkward and confusing to use the touch interaction in front of a screen but not touch it, whereas the other techniques feel more connected to the space above the sensor than the display screen. Getting the touch interaction to engage was difficult, but maintaining the touch once engaged was simple and disengaging felt easy."
$ws.Cells.Item(7, 7).Value = "dwell and pinch anywhere were my favourites, they were intuitive enough that I could focus more on the task and less on performing the gesture, I found dwell easier to use as pinch was hard to get to engage but pinch was better for movement and disengagement "
$ws.Cells.Item(7, 8).Value = "Dwell;Pinch Anywhere;Pinch on Circle;Touch In The Air;"
$ws.Cells.Item(7, 9).Value = "Dwell was most easy to use because it gave a visual indicator that I was engaged with the slider, and once you understand the limit on speed of movement for your dwell to remain engaged it was very easy to use accurately. Disengaging was somewhat difficult however."
$ws.Cells.Item(7, 10).Value = "touch was difficult to use because the sensor did not always register the touch motion, requiring multiple attempts to engage with the slider circle. This meant if it got lost halfway through the task, reengaging with the bar was frustrating"

# Row 8 (participant ID 7)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 44573.742581018501
$ws.Cells.Item(8, 3).Value = 44573.746030092603
$ws.Cells.Item(8, 4).Value = "'7"
$ws.Cells.Item(8, 5).Value = "Touch In The Air;Dwell;Pinch Anywhere;Pinch on Circle;"
$ws.Cells.Item(8, 6).Value = "Was difficult to control, sensor often wasn't responsive to pinch motion which made it difficult to engage and disengage "
$ws.Cells.Item(8, 7).Value = "Easiest to use, benefit of not having to wait to control (like dwell)"
$ws.Cells.Item(8, 8).Value = "Touch In The Air;Dwell;Pinch Anywhere;Pinch on Circle;"
$ws.Cells.Item(8, 9).Value = "Most responsive, most intuitive, instantly able to control so easy to correct small errors"
$ws.Cells.Item(8, 10).Value = "Often unresponsive to pinch"

# --- Apply the same cell styles used by the rest of the table ---
$ws.Range("B7:C8").NumberFormat = "m/d/yy h:mm:ss"

# --- Grow the table / autofilter range to cover the new rows ---
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:J8"))

# --- Update selection to match the saved state ---
$ws.Range("D4").Select()

# --- Update workbook window position ---
$wb.Windows.Item(1).Left = 0
